# Adds a new row to the "Git Basic Commands" cheatsheet table for
# force-pushing to a branch, and expands Table1 to include it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cheatsheet entry: command + description (no "Notes" value).
$ws.Range("A16").Value = "git push --force origin {branch_name}"
$ws.Range("B16").Value = "Force through a commit"

# Grow the worksheet table (ListObject) so the new row becomes part of it.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C16"))

# Match the author's final selection/cursor position.
$ws.Range("C16").Select() | Out-Null
